{"js": "const replacements = [\n  [\"2024-06-17 Monday\", \"2024-06-18 Tuesday\"],\n  [\"518\u00d72=\", \"798\u00d74=\"],\n  [\"736\u00d77=\", \"895\u00d78=\"],\n  [\"902\u00d77=\", \"979\u00d79=\"],\n  [\"962\u00d74=\", \"460\u00d76=\"],\n  [\"563\u00d73=\", \"191\u00d78=\"],\n  [\"995\u00d79=\", \"397\u00d72=\"],\n  [\"979\u00d76=\", \"108\u00d72=\"],\n  [\"142\u00d79=\", \"321\u00d79=\"],\n  [\"827\u00d75=\", \"981\u00d74=\"],\n  [\"274\u00d74=\", \"527\u00d75=\"],\n  [\"199\u00d76=\", \"252\u00d76=\"],\n  [\"481\u00d77=\", \"976\u00d75=\"],\n  [\"958\u00d76=\", \"133\u00d78=\"],\n  [\"347\u00d72=\", \"523\u00d78=\"],\n  [\"907\u00d73=\", \"125\u00d77=\"],\n  [\"990\u00d74=\", \"794\u00d76=\"],\n  [\"805\u00d78=\", \"133\u00d75=\"],\n  [\"253\u00d72=\", \"279\u00d75=\"],\n  [\"709\u00d73=\", \"171\u00d72=\"],\n  [\"602\u00d72=\", \"678\u00d76=\"],\n  [\"120\u00d75=\", \"457\u00d75=\"],\n  [\"838\u00d79=\", \"670\u00d72=\"],\n  [\"618\u00d72=\", \"250\u00d76=\"],\n  [\"111\u00d73=\", \"124\u00d76=\"],\n  [\"317\u00d78=\", \"474\u00d72=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @('2024-06-17 Monday', '2024-06-18 Tuesday'),\n    @('518\u00d72=', '798\u00d74='),\n    @('736\u00d77=', '895\u00d78='),\n    @('902\u00d77=', '979\u00d79='),\n    @('962\u00d74=', '460\u00d76='),\n    @('563\u00d73=', '191\u00d78='),\n    @('995\u00d79=', '397\u00d72='),\n    @('979\u00d76=', '108\u00d72='),\n    @('142\u00d79=', '321\u00d79='),\n    @('827\u00d75=', '981\u00d74='),\n    @('274\u00d74=', '527\u00d75='),\n    @('199\u00d76=', '252\u00d76='),\n    @('481\u00d77=', '976\u00d75='),\n    @('958\u00d76=', '133\u00d78='),\n    @('347\u00d72=', '523\u00d78='),\n    @('907\u00d73=', '125\u00d77='),\n    @('990\u00d74=', '794\u00d76='),\n    @('805\u00d78=', '133\u00d75='),\n    @('253\u00d72=', '279\u00d75='),\n    @('709\u00d73=', '171\u00d72='),\n    @('602\u00d72=', '678\u00d76='),\n    @('120\u00d75=', '457\u00d75='),\n    @('838\u00d79=', '670\u00d72='),\n    @('618\u00d72=', '250\u00d76='),\n    @('111\u00d73=', '124\u00d76='),\n    @('317\u00d78=', '474\u00d72='),\n)\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll)\n    if (-not $found) {\n        throw \"No match found for: $oldText\"\n    }\n}"}
